# Update the Ifnb1-Ifnar1 LR-pairs sheet with refreshed TPM-derived numbers and
# add the "MuSCs" sending-cluster rows alongside the existing "Resolving-Mac" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple is: SendingCluster, LigandSymbol, ReceptorSymbol, TargetCluster,
# then the 16 numeric metric columns (E..T).
$rows = @(
    ,@("MuSCs","Ifnb1","Ifnar1","ECs",1,0.3333333333333333,0.01317833333333333,0.039535,0.2048477437473122,0.2048477437473122,3,1,21.602164,64.80649199999999,0.2504461301095448,0.2504461301095448,0.2846805179133333,2.56212466122,0.05130332468318602,0.05130332468318602)
    ,@("MuSCs","Ifnb1","Ifnar1","FAPs",1,0.3333333333333333,0.01317833333333333,0.039535,0.2048477437473122,0.2048477437473122,3,1,20.520975,61.56292500000001,0.2379112932771326,0.2379112932771326,0.2704322488750001,2.433890239875,0.04873559163982569,0.04873559163982569)
    ,@("MuSCs","Ifnb1","Ifnar1","MuSCs",1,0.3333333333333333,0.01317833333333333,0.039535,0.2048477437473122,0.2048477437473122,3,1,11.57857933333334,34.735738,0.1342370322806403,0.1342370322806404,0.1525863779811111,1.37327740183,0.02749815319002428,0.02749815319002429)
    ,@("MuSCs","Ifnb1","Ifnar1","Resolving-Mac",1,0.3333333333333333,0.01317833333333333,0.039535,0.2048477437473122,0.2048477437473122,3,1,32.55301433333333,97.659043,0.3774055443326823,0.3774055443326824,0.4289944738894444,3.860950265005,0.07731067423427616,0.07731067423427619)
    ,@("Resolving-Mac","Ifnb1","Ifnar1","ECs",2,0.6666666666666666,0.051154,0.153462,0.7951522562526879,0.7951522562526879,3,1,21.602164,64.80649199999999,0.2504461301095448,0.2504461301095448,1.105037097256,9.945333875303998,0.1991428054263588,0.1991428054263588)
    ,@("Resolving-Mac","Ifnb1","Ifnar1","FAPs",2,0.6666666666666666,0.051154,0.153462,0.7951522562526879,0.7951522562526879,3,1,20.520975,61.56292500000001,0.2379112932771326,0.2379112932771326,1.04972995515,9.44756959635,0.1891757016373069,0.189175701637307)
    ,@("Resolving-Mac","Ifnb1","Ifnar1","MuSCs",2,0.6666666666666666,0.051154,0.153462,0.7951522562526879,0.7951522562526879,3,1,11.57857933333334,34.735738,0.1342370322806403,0.1342370322806404,0.5922906472173334,5.330615824956,0.1067388790906161,0.1067388790906161)
    ,@("Resolving-Mac","Ifnb1","Ifnar1","Resolving-Mac",2,0.6666666666666666,0.051154,0.153462,0.7951522562526879,0.7951522562526879,3,1,32.55301433333333,97.659043,0.3774055443326823,0.3774055443326824,1.665216895207333,14.986952056866,0.3000948700984062,0.3000948700984062)
)

$numRows = $rows.Count
$numCols = 20
$data = New-Object 'object[,]' $numRows,$numCols
for ($i = 0; $i -lt $numRows; $i++) {
    for ($j = 0; $j -lt $numCols; $j++) {
        $data[$i,$j] = $rows[$i][$j]
    }
}

# Data now spans rows 2-9 (header stays in row 1); dimension grows from T5 to T9.
$ws.Range("A2:T9").Value = $data
